$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) cells whose new values look numeric,
# so Excel keeps them as literal text (matching the original formatting),
# then set the Coin/Link/Price/Volume values for each changed row.

$ws.Range("D2").Value = '29.361.83'
$ws.Range("E2").Value = '  +1.65%  '
$ws.Range("D3").Value = '1.833.66'
$ws.Range("E3").Value = '  +0.69%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9981'
$ws.Range("E4").Value = '  +0.57%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.33'
$ws.Range("E5").Value = '  +1.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6212'
$ws.Range("E6").Value = '  +0.93%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9993'
$ws.Range("E7").Value = '  +0.42%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07385'
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2911'
$ws.Range("E9").Value = '  -0.54%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.34'
$ws.Range("E10").Value = '  +1.94%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07657'
$ws.Range("E11").Value = '  +0.54%  '
$ws.Range("D12").Value = '1.846.26'
$ws.Range("E12").Value = '  +0.80%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.003'
$ws.Range("E13").Value = '  +0.58%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6761'
$ws.Range("E14").Value = '  +0.86%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '82.95'
$ws.Range("E15").Value = '  +0.69%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000008990'
$ws.Range("E16").Value = '  -0.63%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.904'
$ws.Range("E17").Value = '  +0.87%  '
$ws.Range("D18").Value = '29.339.03'
$ws.Range("E18").Value = '  +1.41%  '
$ws.Range("D19").Value = '2.090.53'
$ws.Range("E19").Value = '  -0.24%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '241.66'
$ws.Range("E20").Value = '  +1.53%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.54'
$ws.Range("E21").Value = '  -0.70%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9992'
$ws.Range("E22").Value = '  +0.42%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.387'
$ws.Range("E23").Value = '  +2.88%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9989'
$ws.Range("E24").Value = '  +0.55%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.33'
$ws.Range("E25").Value = '  +0.25%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.578'
$ws.Range("E26").Value = '  +1.58%  '
$ws.Range("B27").Value = 'Stellar'
$ws.Range("C27").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1396'
$ws.Range("E27").Value = '  -0.68%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.73'
$ws.Range("E28").Value = '  -0.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.491'
$ws.Range("E29").Value = '  +0.36%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05926'
$ws.Range("E30").Value = '  +6.81%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.094'
$ws.Range("E31").Value = '  +0.20%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.111'
$ws.Range("E32").Value = '  +0.12%  '
$ws.Range("B33").Value = 'Toncoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.224'
$ws.Range("E33").Value = '  +1.94%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.868'
$ws.Range("E34").Value = '  +2.15%  '
$ws.Range("E35").Value = '  +0.71%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7215'
$ws.Range("E36").Value = '  -1.91%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.615'
$ws.Range("E37").Value = '  -0.65%  '
$ws.Range("E38").Value = '  +3.87%  '
$ws.Range("D39").Value = '1.225.77'
$ws.Range("E39").Value = '  +2.05%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01769'
$ws.Range("E40").Value = '  -0.07%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9161'
$ws.Range("E41").Value = '  +2.76%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.221'
$ws.Range("E42").Value = '  -2.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.000'
$ws.Range("E43").Value = '  +0.67%  '
$ws.Range("D44").Value = '2.002.63'
$ws.Range("E44").Value = '  +0.68%  '
$ws.Range("E45").Value = '  +1.17%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '65.75'
$ws.Range("E46").Value = '  +1.04%  '
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000120'
$ws.Range("E47").Value = '  -1.56%  '
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5054'
$ws.Range("E48").Value = '  -0.04%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.232'
$ws.Range("E49").Value = '  +1.87%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4058'
$ws.Range("E50").Value = '  +0.57%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1173'
$ws.Range("E51").Value = '  +6.74%  '
